# Update "Australia ALeague" worksheet data:
#  - delete the last data row (row 149, the Melbourne City vs Perth Glory
#    fixture that is duplicated further up after the refresh)
#  - refresh the odds/results for rows 146-148 with the newer data pull

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete last row (id 147 / match 7127409 old copy).
$ws.Rows(149).Delete()

# --- Row 146: match 7127404 (Wellington Phoenix vs Melbourne Victory) ---
$ws.Cells.Item(146, 2).Value  = 7127404
$ws.Cells.Item(146, 5).Value  = 45394.16666666666
$ws.Cells.Item(146, 6).Value  = "Wellington Phoenix"
$ws.Cells.Item(146, 7).Value  = "Melbourne Victory"
$ws.Cells.Item(146, 8).Value  = 1
$ws.Cells.Item(146, 9).Value  = 0
$ws.Cells.Item(146, 10).Value = "H"
$ws.Cells.Item(146, 11).Value = 3.1
$ws.Cells.Item(146, 12).Value = 3.4
$ws.Cells.Item(146, 13).Value = 2.25
$ws.Cells.Item(146, 14).Value = 2.9
$ws.Cells.Item(146, 15).Value = 3.6
$ws.Cells.Item(146, 16).Value = 2.3
$ws.Cells.Item(146, 17).Value = 0.25
$ws.Cells.Item(146, 18).Value = 1.86
$ws.Cells.Item(146, 19).Value = 2.04
$ws.Cells.Item(146, 20).Value = 2.75
$ws.Cells.Item(146, 21).Value = 1.925
$ws.Cells.Item(146, 22).Value = 1.925
$ws.Cells.Item(146, 23).Value = 1.9
$ws.Cells.Item(146, 24).Value = -1
$ws.Cells.Item(146, 25).Value = -1
$ws.Cells.Item(146, 26).Value = 0.8600000000000001
$ws.Cells.Item(146, 27).Value = -1
$ws.Cells.Item(146, 28).Value = -1
$ws.Cells.Item(146, 29).Value = 0.925

# --- Row 147: match 7127405 (Adelaide United vs Macarthur FC) ---
$ws.Cells.Item(147, 2).Value  = 7127405
$ws.Cells.Item(147, 5).Value  = 45394.28125
$ws.Cells.Item(147, 6).Value  = "Adelaide United"
$ws.Cells.Item(147, 7).Value  = "Macarthur FC"
$ws.Cells.Item(147, 8).Value  = 1
$ws.Cells.Item(147, 9).Value  = 2
$ws.Cells.Item(147, 10).Value = "A"
$ws.Cells.Item(147, 11).Value = 1.833
$ws.Cells.Item(147, 12).Value = 4
$ws.Cells.Item(147, 13).Value = 3.75
$ws.Cells.Item(147, 14).Value = 1.5
$ws.Cells.Item(147, 15).Value = 5.25
$ws.Cells.Item(147, 16).Value = 5
$ws.Cells.Item(147, 17).Value = -1.25
$ws.Cells.Item(147, 18).Value = 2
$ws.Cells.Item(147, 19).Value = 1.85
$ws.Cells.Item(147, 20).Value = 4
$ws.Cells.Item(147, 21).Value = 1.925
$ws.Cells.Item(147, 22).Value = 1.925
$ws.Cells.Item(147, 23).Value = -1
$ws.Cells.Item(147, 24).Value = -1
$ws.Cells.Item(147, 25).Value = 4
$ws.Cells.Item(147, 26).Value = -1
$ws.Cells.Item(147, 27).Value = 0.8500000000000001
$ws.Cells.Item(147, 28).Value = -1
$ws.Cells.Item(147, 29).Value = 0.925

# --- Row 148: match 7127409 (Melbourne City vs Perth Glory) ---
$ws.Cells.Item(148, 2).Value  = 7127409
$ws.Cells.Item(148, 5).Value  = 45396.08333333334
$ws.Cells.Item(148, 6).Value  = "Melbourne City"
$ws.Cells.Item(148, 7).Value  = "Perth Glory"
$ws.Cells.Item(148, 11).Value = 1.571
$ws.Cells.Item(148, 12).Value = 4.5
$ws.Cells.Item(148, 13).Value = 4.75
$ws.Cells.Item(148, 14).Value = 1.4
$ws.Cells.Item(148, 15).Value = 5
$ws.Cells.Item(148, 16).Value = 7
$ws.Cells.Item(148, 17).Value = -1.5
$ws.Cells.Item(148, 18).Value = 2.01
$ws.Cells.Item(148, 19).Value = 1.89
$ws.Cells.Item(148, 20).Value = 3.5
$ws.Cells.Item(148, 21).Value = 1.85
$ws.Cells.Item(148, 22).Value = 2
